# Saldo_guide.xlsx refresh
# - Rename sheet to reflect the new extraction run (20240930-091049)
# - Refresh "Dt. Referencia" (column G) on every data row to the new reference date
# - Update "Saldo Previsto" / "Vl. Total" (columns E and H) for accounts whose
#   balances moved between the 2024-09-27 and 2024-09-30 extraction runs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new export timestamp
$ws.Name = "IClientBalance-20240930-091049-"

# New reference date (Excel serial date number for 2024-09-30) applied to every data row
$newDate = 45565
$firstRow = 2
$lastRow = 274
$dateRange = "G" + $firstRow + ":G" + $lastRow
$ws.Range($dateRange).Value2 = $newDate

# Accounts with updated balances (column E = Saldo Previsto, column H = Vl. Total)
$balanceUpdates = @{
    6   = 1059.1199999999999
    15  = 1570.01
    17  = 1174.81
    43  = 813.03
    49  = 879.55
    51  = 1017.12
    52  = 915.19
    57  = 3970.74
    97  = 1256.7
    101 = 3884.66
    102 = 773.16
    105 = 642.74
    107 = 1177.8599999999999
    108 = 403.89
    109 = 324.19
    110 = 820.52
    120 = 972.46
    138 = 32313.71
    143 = 1901.96
    230 = 847.86
    255 = 27474.84
}

foreach ($row in $balanceUpdates.Keys) {
    $value = $balanceUpdates[$row]
    $ws.Range("E$row").Value2 = $value
    $ws.Range("H$row").Value2 = $value
}
